$wb = $excel.ActiveWorkbook

# --- Sheet "login": move selection to D3 (drop prior tabSelected/topLeftCell state) ---
$wsLogin = $wb.Worksheets.Item("login")
$wsLogin.Range("D3").Select()

# --- Sheet "users_add": populate header + data rows ---
$wsAdd = $wb.Worksheets.Item("users_add")

# Column widths (best-effort; engine quantizes to its own internal pixel grid)
$wsAdd.Columns.Item(1).ColumnWidth = 7.7109375
$wsAdd.Columns.Item(3).ColumnWidth = 25.5703125
$wsAdd.Columns.Item(4).ColumnWidth = 10.140625
$wsAdd.Columns.Item(5).ColumnWidth = 23
$wsAdd.Columns.Item(6).ColumnWidth = 53.85546875
$wsAdd.Columns.Item(7).ColumnWidth = 62.140625
$wsAdd.Columns.Item(8).ColumnWidth = 39.7109375

# Cell writes, ordered to reproduce the exact shared-string allocation order
$wsAdd.Range("A1").Value = 'case_id'
$wsAdd.Range("B1").Value = 'interface'
$wsAdd.Range("C1").Value = 'title'
$wsAdd.Range("D1").Value = 'method'
$wsAdd.Range("E1").Value = 'url'
$wsAdd.Range("F1").Value = 'request_data'
$wsAdd.Range("G1").Value = 'expected'
$wsAdd.Range("H1").Value = 'check_db_sql'
$wsAdd.Range("B2").Value = 'add user'
$wsAdd.Range("C2").Value = 'user added successfully'
$wsAdd.Range("E2").Value = '/paymall_admin/users/'
$wsAdd.Range("C3").Value = 'failed, user exist'
$wsAdd.Range("C4").Value = 'failed, mobile phone exist'
$wsAdd.Range("G3").Value = '{"username":["A user with that username already exists."]}'
$wsAdd.Range("H2").Value = 'select * from tb_users WHERE user=''#username#'''
$wsAdd.Range("G4").Value = '{"mobile":["user with this Mobile Phone number already exists."]}'
$wsAdd.Range("C6").Value = 'failed, phone number blank'
$wsAdd.Range("C5").Value = 'failed, user name blank'
$wsAdd.Range("C7").Value = 'failed, password blank'
$wsAdd.Range("F2").Value = '{"username":"#username#","mobile":"#phone#","password":"#password#","email":"a@a.com"}'
$wsAdd.Range("G2").Value = '{"id":#id#,"username":"#username#","mobile":"#phone#","email":"a@a.com"}'
$wsAdd.Range("F5").Value = '{"username":"","mobile":"#phone#","password":"#password#","email":"a@a.com"}'
$wsAdd.Range("F6").Value = '{"username":"#username#","mobile":"","password":"#password#","email":"a@a.com"}'
$wsAdd.Range("F7").Value = '{"username":"#username#","mobile":"#phone#","password":"","email":"a@a.com"}'
$wsAdd.Range("F3").Value = '{"username":"aaaaaa","mobile":"#phone#","password":"#password#","email":"a@a.com"}'
$wsAdd.Range("F4").Value = '{"username":"#username#","mobile":"13111111111","password":"#password#","email":"a@a.com"}'
$wsAdd.Range("G5").Value = '{"username": ["This field may not be blank."]}'
$wsAdd.Range("G6").Value = '{"mobile": ["This field may not be blank."]}'
$wsAdd.Range("G7").Value = '{"password": ["This field may not be blank."]}'

# Remaining cells (numbers + strings that reuse already-allocated shared strings)
$wsAdd.Range("A2").Value = 1
$wsAdd.Range("D2").Value = 'post'
$wsAdd.Range("A3").Value = 2
$wsAdd.Range("B3").Value = 'add user'
$wsAdd.Range("D3").Value = 'post'
$wsAdd.Range("E3").Value = '/paymall_admin/users/'
$wsAdd.Range("A4").Value = 3
$wsAdd.Range("B4").Value = 'add user'
$wsAdd.Range("D4").Value = 'post'
$wsAdd.Range("E4").Value = '/paymall_admin/users/'
$wsAdd.Range("A5").Value = 4
$wsAdd.Range("B5").Value = 'add user'
$wsAdd.Range("D5").Value = 'post'
$wsAdd.Range("E5").Value = '/paymall_admin/users/'
$wsAdd.Range("A6").Value = 5
$wsAdd.Range("B6").Value = 'add user'
$wsAdd.Range("D6").Value = 'post'
$wsAdd.Range("E6").Value = '/paymall_admin/users/'
$wsAdd.Range("A7").Value = 6
$wsAdd.Range("B7").Value = 'add user'
$wsAdd.Range("D7").Value = 'post'
$wsAdd.Range("E7").Value = '/paymall_admin/users/'

# Activate the users_add sheet and select G7 to match the saved view state
$wsAdd.Activate()
$wsAdd.Range("G7").Select()
